{"js": "const doc = context.document;\nconst body = doc.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\n// Helper: find the (single) paragraph whose text contains `needle` and\n// append `addition` to it as a brand-new run (InsertLocation.end keeps\n// the existing runs untouched and adds a separate trailing run, matching\n// how the document was actually edited).\nfunction findParagraph(needle) {\n  const matches = paragraphs.items.filter((p) => p.text.indexOf(needle) !== -1);\n  if (matches.length !== 1) {\n    throw new Error(\n      \"Expected exactly one paragraph containing \" + JSON.stringify(needle) +\n      \" but found \" + matches.length\n    );\n  }\n  return matches[0];\n}\n\nfindParagraph(\"t=1, 9, 17, 25, 33 - every 8\").insertText(\" (all odd)\", Word.InsertLocation.end);\nfindParagraph(\"f=8, 10, 16, 18, 24\").insertText(\"(all even)\", Word.InsertLocation.end);\nfindParagraph(\"m=3, 7, 11, 15, 19\").insertText(\" (all even)\", Word.InsertLocation.end);\nfindParagraph(\"r = 4, 6, 12, 14, 20\").insertText(\"v(all even)\", Word.InsertLocation.end);\nconst lastLine = findParagraph(\"l=5, 13, 21, 29, 37\");\nlastLine.insertText(\"(all odd)\", Word.InsertLocation.end);\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark from its old position (mid-sentence in the\n// \"In coming up with a pattern...\" paragraph) to the end of the line we\n// just annotated above.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst endOfLastLine = lastLine.getRange(Word.RangeLocation.end);\nendOfLastLine.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Append-AfterText($needle, $newText) {\n    $rng = $d.Content\n    $rng.Find.Text = $needle\n    $found = $rng.Find.Execute()\n    if (-not $found) {\n        throw \"Text not found: $needle\"\n    }\n    $rng.Collapse(0)\n    $rng.InsertAfter($newText)\n    return $rng\n}\n\n# Append the \"(all odd)\"/\"(all even)\" annotations to each of the five\n# finger-pattern lines.\nAppend-AfterText \"t=1, 9, 17, 25, 33 - every 8\" \" (all odd)\"\nAppend-AfterText \"f=8, 10, 16, 18, 24 \u2013 every 2 & every 4\" \"(all even)\"\nAppend-AfterText \"m=3, 7, 11, 15, 19 \u2013 every 4\" \" (all even)\"\nAppend-AfterText \"r = 4, 6, 12, 14, 20 \u2013 every 2 & every 4\" \"v(all even)\"\n$lastRng = Append-AfterText \"l=5, 13, 21, 29, 37 \u2013 every 8\" \"(all odd)\"\n\n# Move the \"_GoBack\" bookmark from its old spot (mid-sentence in the\n# paragraph that starts \"In coming up with a pattern...\") to the end of\n# the line we just annotated (\"l=5, ... (all odd)\"), as a collapsed\n# (zero-length) bookmark.\n#\n# A collapsed Range sitting exactly on a paragraph-mark offset cannot be\n# used directly to seed Bookmarks.Add, so a one-character placeholder is\n# inserted, bookmarked, and then removed -- leaving the bookmark collapsed\n# at the desired position.\n$lastRng.InsertAfter(\"X\")\n$placeholder = $d.Range($lastRng.End - 1, $lastRng.End)\n$d.Bookmarks.Add(\"_GoBack\", $placeholder)\n$placeholder.Text = \"\"\n"}
